$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn on page setup (paper size A4 / portrait) to mirror the authored workbook change.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# The sheet is driven by an XML-mapped table ("表1"); grow it by two columns
# (new NPC fields: CanClone, ActorID) using the ListObject/table object model.
$lo = $ws.ListObjects.Item(1)
$colCanClone = $lo.ListColumns.Add()
$colActorID = $lo.ListColumns.Add()

# Header text - set K1 ("ActorID") before J1 ("CanClone") so the shared-string
# table receives the two new labels in that same order.
$ws.Range("K1").Value = "ActorID"
$ws.Range("J1").Value = "CanClone"

# Data rows for the two new columns.
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 0

$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0

# Column J gets an explicit (non bestFit) width of 14 characters.
$ws.Columns.Item(10).ColumnWidth = 13.285714285714286

# Match the author's final selection/cursor position.
$ws.Range("K9").Select() | Out-Null
